$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 7: name in A7, URL text (with hyperlink) in B7
$ws.Range("A7").Value = "Grand Piano"
$ws.Range("B7").Value = "https://www.iconshock.com/musical-instruments-icons/"

$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.iconshock.com/musical-instruments-icons/") | Out-Null

# Ensure B7 uses the same "Hyperlink" cell style as the other link cells (B1:B6)
$ws.Range("B7").Style = "Hyperlink"
